# Update scripts with new TPM-derived NATMI values for Spon2-Itga4.
# Re-computed ligand/receptor expression stats (columns E:T) for every
# Sending-cluster x Target-cluster pair, and add the new "Inflammatory-Mac"
# sending-cluster block (rows 18-21) that the refreshed TPM table now
# produces.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 1).Value = "ECs"
$ws.Cells.Item(2, 2).Value = "Spon2"
$ws.Cells.Item(2, 3).Value = "Itga4"
$ws.Cells.Item(2, 4).Value = "ECs"
$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 6).Value = 1
$ws.Cells.Item(2, 7).Value = 0.54332
$ws.Cells.Item(2, 8).Value = 1.62996
$ws.Cells.Item(2, 9).Value = 0.09277509850694737
$ws.Cells.Item(2, 10).Value = 0.09480543614915297
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 12).Value = 1
$ws.Cells.Item(2, 13).Value = 0.506715
$ws.Cells.Item(2, 14).Value = 1.520145
$ws.Cells.Item(2, 15).Value = 0.003122343715987576
$ws.Cells.Item(2, 16).Value = 0.003132472094339857
$ws.Cells.Item(2, 17).Value = 0.2753083938000001
$ws.Cells.Item(2, 18).Value = 2.4777755442
$ws.Cells.Item(2, 19).Value = 0.0002896757458232955
$ws.Cells.Item(2, 20).Value = 0.0002969753831289408

$ws.Cells.Item(3, 1).Value = "ECs"
$ws.Cells.Item(3, 2).Value = "Spon2"
$ws.Cells.Item(3, 3).Value = "Itga4"
$ws.Cells.Item(3, 4).Value = "Inflammatory-Mac"
$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 6).Value = 1
$ws.Cells.Item(3, 7).Value = 0.54332
$ws.Cells.Item(3, 8).Value = 1.62996
$ws.Cells.Item(3, 9).Value = 0.09277509850694737
$ws.Cells.Item(3, 10).Value = 0.09480543614915297
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(3, 13).Value = 88.13219433333332
$ws.Cells.Item(3, 14).Value = 264.396583
$ws.Cells.Item(3, 15).Value = 0.5430646480820168
$ws.Cells.Item(3, 16).Value = 0.5448262620252092
$ws.Cells.Item(3, 17).Value = 47.88398382518666
$ws.Cells.Item(3, 18).Value = 430.9558544266799
$ws.Cells.Item(3, 19).Value = 0.05038287622144982
$ws.Cells.Item(3, 20).Value = 0.05165249139681266

$ws.Cells.Item(4, 1).Value = "ECs"
$ws.Cells.Item(4, 2).Value = "Spon2"
$ws.Cells.Item(4, 3).Value = "Itga4"
$ws.Cells.Item(4, 4).Value = "MuSCs"
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 6).Value = 1
$ws.Cells.Item(4, 7).Value = 0.54332
$ws.Cells.Item(4, 8).Value = 1.62996
$ws.Cells.Item(4, 9).Value = 0.09277509850694737
$ws.Cells.Item(4, 10).Value = 0.09480543614915297
$ws.Cells.Item(4, 11).Value = 2
$ws.Cells.Item(4, 12).Value = 1
$ws.Cells.Item(4, 13).Value = 1.5741895
$ws.Cells.Item(4, 14).Value = 3.148379
$ws.Cells.Item(4, 15).Value = 0.009700049718478087
$ws.Cells.Item(4, 16).Value = 0.006487676741301404
$ws.Cells.Item(4, 17).Value = 0.85528863914
$ws.Cells.Item(4, 18).Value = 5.13173183484
$ws.Cells.Item(4, 19).Value = 0.0008999230681540916
$ws.Cells.Item(4, 20).Value = 0.0006150670230537952

$ws.Cells.Item(5, 1).Value = "ECs"
$ws.Cells.Item(5, 2).Value = "Spon2"
$ws.Cells.Item(5, 3).Value = "Itga4"
$ws.Cells.Item(5, 4).Value = "Resolving-Mac"
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 0.54332
$ws.Cells.Item(5, 8).Value = 1.62996
$ws.Cells.Item(5, 9).Value = 0.09277509850694737
$ws.Cells.Item(5, 10).Value = 0.09480543614915297
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 12).Value = 1
$ws.Cells.Item(5, 13).Value = 72.07364666666666
$ws.Cells.Item(5, 14).Value = 216.22094
$ws.Cells.Item(5, 15).Value = 0.4441129584835175
$ws.Cells.Item(5, 16).Value = 0.4455535891391496
$ws.Cells.Item(5, 17).Value = 39.15905370693334
$ws.Cells.Item(5, 18).Value = 352.4314833624
$ws.Cells.Item(5, 19).Value = 0.04120262347152016
$ws.Cells.Item(5, 20).Value = 0.04224090234615758

$ws.Cells.Item(6, 1).Value = "FAPs"
$ws.Cells.Item(6, 2).Value = "Spon2"
$ws.Cells.Item(6, 3).Value = "Itga4"
$ws.Cells.Item(6, 4).Value = "ECs"
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 6).Value = 1
$ws.Cells.Item(6, 7).Value = 4.616901666666666
$ws.Cells.Item(6, 8).Value = 13.850705
$ws.Cells.Item(6, 9).Value = 0.7883632241071366
$ws.Cells.Item(6, 10).Value = 0.805616167573593
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 12).Value = 1
$ws.Cells.Item(6, 13).Value = 0.506715
$ws.Cells.Item(6, 14).Value = 1.520145
$ws.Cells.Item(6, 15).Value = 0.003122343715987576
$ws.Cells.Item(6, 16).Value = 0.003132472094339857
$ws.Cells.Item(6, 17).Value = 2.339453328025
$ws.Cells.Item(6, 18).Value = 21.055079952225
$ws.Cells.Item(6, 19).Value = 0.002461540958706623
$ws.Cells.Item(6, 20).Value = 0.002523570163673302

$ws.Cells.Item(7, 1).Value = "FAPs"
$ws.Cells.Item(7, 2).Value = "Spon2"
$ws.Cells.Item(7, 3).Value = "Itga4"
$ws.Cells.Item(7, 4).Value = "Inflammatory-Mac"
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 6).Value = 1
$ws.Cells.Item(7, 7).Value = 4.616901666666666
$ws.Cells.Item(7, 8).Value = 13.850705
$ws.Cells.Item(7, 9).Value = 0.7883632241071366
$ws.Cells.Item(7, 10).Value = 0.805616167573593
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 12).Value = 1
$ws.Cells.Item(7, 13).Value = 88.13219433333332
$ws.Cells.Item(7, 14).Value = 264.396583
$ws.Cells.Item(7, 15).Value = 0.5430646480820168
$ws.Cells.Item(7, 16).Value = 0.5448262620252092
$ws.Cells.Item(7, 17).Value = 406.8976749045571
$ws.Cells.Item(7, 18).Value = 3662.079074141014
$ws.Cells.Item(7, 19).Value = 0.4281321968605463
$ws.Cells.Item(7, 20).Value = 0.4389208452061952

$ws.Cells.Item(8, 1).Value = "FAPs"
$ws.Cells.Item(8, 2).Value = "Spon2"
$ws.Cells.Item(8, 3).Value = "Itga4"
$ws.Cells.Item(8, 4).Value = "MuSCs"
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 6).Value = 1
$ws.Cells.Item(8, 7).Value = 4.616901666666666
$ws.Cells.Item(8, 8).Value = 13.850705
$ws.Cells.Item(8, 9).Value = 0.7883632241071366
$ws.Cells.Item(8, 10).Value = 0.805616167573593
$ws.Cells.Item(8, 11).Value = 2
$ws.Cells.Item(8, 12).Value = 1
$ws.Cells.Item(8, 13).Value = 1.5741895
$ws.Cells.Item(8, 14).Value = 3.148379
$ws.Cells.Item(8, 15).Value = 0.009700049718478087
$ws.Cells.Item(8, 16).Value = 0.006487676741301404
$ws.Cells.Item(8, 17).Value = 7.267878126199165
$ws.Cells.Item(8, 18).Value = 43.60726875719499
$ws.Cells.Item(8, 19).Value = 0.007647162470058907
$ws.Cells.Item(8, 20).Value = 0.005226577272783574

$ws.Cells.Item(9, 1).Value = "FAPs"
$ws.Cells.Item(9, 2).Value = "Spon2"
$ws.Cells.Item(9, 3).Value = "Itga4"
$ws.Cells.Item(9, 4).Value = "Resolving-Mac"
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 6).Value = 1
$ws.Cells.Item(9, 7).Value = 4.616901666666666
$ws.Cells.Item(9, 8).Value = 13.850705
$ws.Cells.Item(9, 9).Value = 0.7883632241071366
$ws.Cells.Item(9, 10).Value = 0.805616167573593
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 12).Value = 1
$ws.Cells.Item(9, 13).Value = 72.07364666666666
$ws.Cells.Item(9, 14).Value = 216.22094
$ws.Cells.Item(9, 15).Value = 0.4441129584835175
$ws.Cells.Item(9, 16).Value = 0.4455535891391496
$ws.Cells.Item(9, 17).Value = 332.7569394180777
$ws.Cells.Item(9, 18).Value = 2994.8124547627
$ws.Cells.Item(9, 19).Value = 0.3501223238178247
$ws.Cells.Item(9, 20).Value = 0.3589451749309409

$ws.Cells.Item(10, 1).Value = "Inflammatory-Mac"
$ws.Cells.Item(10, 2).Value = "Spon2"
$ws.Cells.Item(10, 3).Value = "Itga4"
$ws.Cells.Item(10, 4).Value = "ECs"
$ws.Cells.Item(10, 5).Value = 1
$ws.Cells.Item(10, 6).Value = 0.3333333333333333
$ws.Cells.Item(10, 7).Value = 0.073119
$ws.Cells.Item(10, 8).Value = 0.219357
$ws.Cells.Item(10, 9).Value = 0.01248550104492653
$ws.Cells.Item(10, 10).Value = 0.01275874012697842
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 12).Value = 1
$ws.Cells.Item(10, 13).Value = 0.506715
$ws.Cells.Item(10, 14).Value = 1.520145
$ws.Cells.Item(10, 15).Value = 0.003122343715987576
$ws.Cells.Item(10, 16).Value = 0.003132472094339857
$ws.Cells.Item(10, 17).Value = 0.037050494085
$ws.Cells.Item(10, 18).Value = 0.333454446765
$ws.Cells.Item(10, 19).Value = 0.00003898402572858268
$ws.Cells.Item(10, 20).Value = 0.00003996639740669406

$ws.Cells.Item(11, 1).Value = "Inflammatory-Mac"
$ws.Cells.Item(11, 2).Value = "Spon2"
$ws.Cells.Item(11, 3).Value = "Itga4"
$ws.Cells.Item(11, 4).Value = "Inflammatory-Mac"
$ws.Cells.Item(11, 5).Value = 1
$ws.Cells.Item(11, 6).Value = 0.3333333333333333
$ws.Cells.Item(11, 7).Value = 0.073119
$ws.Cells.Item(11, 8).Value = 0.219357
$ws.Cells.Item(11, 9).Value = 0.01248550104492653
$ws.Cells.Item(11, 10).Value = 0.01275874012697842
$ws.Cells.Item(11, 11).Value = 3
$ws.Cells.Item(11, 12).Value = 1
$ws.Cells.Item(11, 13).Value = 88.13219433333332
$ws.Cells.Item(11, 14).Value = 264.396583
$ws.Cells.Item(11, 15).Value = 0.5430646480820168
$ws.Cells.Item(11, 16).Value = 0.5448262620252092
$ws.Cells.Item(11, 17).Value = 6.444137917458999
$ws.Cells.Item(11, 18).Value = 57.99724125713099
$ws.Cells.Item(11, 19).Value = 0.006780434231090681
$ws.Cells.Item(11, 20).Value = 0.006951296691532696

$ws.Cells.Item(12, 1).Value = "Inflammatory-Mac"
$ws.Cells.Item(12, 2).Value = "Spon2"
$ws.Cells.Item(12, 3).Value = "Itga4"
$ws.Cells.Item(12, 4).Value = "MuSCs"
$ws.Cells.Item(12, 5).Value = 1
$ws.Cells.Item(12, 6).Value = 0.3333333333333333
$ws.Cells.Item(12, 7).Value = 0.073119
$ws.Cells.Item(12, 8).Value = 0.219357
$ws.Cells.Item(12, 9).Value = 0.01248550104492653
$ws.Cells.Item(12, 10).Value = 0.01275874012697842
$ws.Cells.Item(12, 11).Value = 2
$ws.Cells.Item(12, 12).Value = 1
$ws.Cells.Item(12, 13).Value = 1.5741895
$ws.Cells.Item(12, 14).Value = 3.148379
$ws.Cells.Item(12, 15).Value = 0.009700049718478087
$ws.Cells.Item(12, 16).Value = 0.006487676741301404
$ws.Cells.Item(12, 17).Value = 0.1151031620505
$ws.Cells.Item(12, 18).Value = 0.690618972303
$ws.Cells.Item(12, 19).Value = 0.0001211099808958975
$ws.Cells.Item(12, 20).Value = 0.00008277458157010683

$ws.Cells.Item(13, 1).Value = "Inflammatory-Mac"
$ws.Cells.Item(13, 2).Value = "Spon2"
$ws.Cells.Item(13, 3).Value = "Itga4"
$ws.Cells.Item(13, 4).Value = "Resolving-Mac"
$ws.Cells.Item(13, 5).Value = 1
$ws.Cells.Item(13, 6).Value = 0.3333333333333333
$ws.Cells.Item(13, 7).Value = 0.073119
$ws.Cells.Item(13, 8).Value = 0.219357
$ws.Cells.Item(13, 9).Value = 0.01248550104492653
$ws.Cells.Item(13, 10).Value = 0.01275874012697842
$ws.Cells.Item(13, 11).Value = 3
$ws.Cells.Item(13, 12).Value = 1
$ws.Cells.Item(13, 13).Value = 72.07364666666666
$ws.Cells.Item(13, 14).Value = 216.22094
$ws.Cells.Item(13, 15).Value = 0.4441129584835175
$ws.Cells.Item(13, 16).Value = 0.4455535891391496
$ws.Cells.Item(13, 17).Value = 5.269952970619999
$ws.Cells.Item(13, 18).Value = 47.42957673557999
$ws.Cells.Item(13, 19).Value = 0.005544972807211372
$ws.Cells.Item(13, 20).Value = 0.005684702456468924

$ws.Cells.Item(14, 1).Value = "MuSCs"
$ws.Cells.Item(14, 2).Value = "Spon2"
$ws.Cells.Item(14, 3).Value = "Itga4"
$ws.Cells.Item(14, 4).Value = "ECs"
$ws.Cells.Item(14, 5).Value = 2
$ws.Cells.Item(14, 6).Value = 1
$ws.Cells.Item(14, 7).Value = 0.3762535
$ws.Cells.Item(14, 8).Value = 0.752507
$ws.Cells.Item(14, 9).Value = 0.0642475070420447
$ws.Cells.Item(14, 10).Value = 0.0437690215344491
$ws.Cells.Item(14, 11).Value = 3
$ws.Cells.Item(14, 12).Value = 1
$ws.Cells.Item(14, 13).Value = 0.506715
$ws.Cells.Item(14, 14).Value = 1.520145
$ws.Cells.Item(14, 15).Value = 0.003122343715987576
$ws.Cells.Item(14, 16).Value = 0.003132472094339857
$ws.Cells.Item(14, 17).Value = 0.1906532922525
$ws.Cells.Item(14, 18).Value = 1.143919753515
$ws.Cells.Item(14, 19).Value = 0.0002006027998805958
$ws.Cells.Item(14, 20).Value = 0.0001371052385532221

$ws.Cells.Item(15, 1).Value = "MuSCs"
$ws.Cells.Item(15, 2).Value = "Spon2"
$ws.Cells.Item(15, 3).Value = "Itga4"
$ws.Cells.Item(15, 4).Value = "Inflammatory-Mac"
$ws.Cells.Item(15, 5).Value = 2
$ws.Cells.Item(15, 6).Value = 1
$ws.Cells.Item(15, 7).Value = 0.3762535
$ws.Cells.Item(15, 8).Value = 0.752507
$ws.Cells.Item(15, 9).Value = 0.0642475070420447
$ws.Cells.Item(15, 10).Value = 0.0437690215344491
$ws.Cells.Item(15, 11).Value = 3
$ws.Cells.Item(15, 12).Value = 1
$ws.Cells.Item(15, 13).Value = 88.13219433333332
$ws.Cells.Item(15, 14).Value = 264.396583
$ws.Cells.Item(15, 15).Value = 0.5430646480820168
$ws.Cells.Item(15, 16).Value = 0.5448262620252092
$ws.Cells.Item(15, 17).Value = 33.16004658059683
$ws.Cells.Item(15, 18).Value = 198.960279483581
$ws.Cells.Item(15, 19).Value = 0.0348905498019349
$ws.Cells.Item(15, 20).Value = 0.02384651239511479

$ws.Cells.Item(16, 1).Value = "MuSCs"
$ws.Cells.Item(16, 2).Value = "Spon2"
$ws.Cells.Item(16, 3).Value = "Itga4"
$ws.Cells.Item(16, 4).Value = "MuSCs"
$ws.Cells.Item(16, 5).Value = 2
$ws.Cells.Item(16, 6).Value = 1
$ws.Cells.Item(16, 7).Value = 0.3762535
$ws.Cells.Item(16, 8).Value = 0.752507
$ws.Cells.Item(16, 9).Value = 0.0642475070420447
$ws.Cells.Item(16, 10).Value = 0.0437690215344491
$ws.Cells.Item(16, 11).Value = 2
$ws.Cells.Item(16, 12).Value = 1
$ws.Cells.Item(16, 13).Value = 1.5741895
$ws.Cells.Item(16, 14).Value = 3.148379
$ws.Cells.Item(16, 15).Value = 0.009700049718478087
$ws.Cells.Item(16, 16).Value = 0.006487676741301404
$ws.Cells.Item(16, 17).Value = 0.59229430903825
$ws.Cells.Item(16, 18).Value = 2.369177236153
$ws.Cells.Item(16, 19).Value = 0.0006232040125961045
$ws.Cells.Item(16, 20).Value = 0.0002839592629985657

$ws.Cells.Item(17, 1).Value = "MuSCs"
$ws.Cells.Item(17, 2).Value = "Spon2"
$ws.Cells.Item(17, 3).Value = "Itga4"
$ws.Cells.Item(17, 4).Value = "Resolving-Mac"
$ws.Cells.Item(17, 5).Value = 2
$ws.Cells.Item(17, 6).Value = 1
$ws.Cells.Item(17, 7).Value = 0.3762535
$ws.Cells.Item(17, 8).Value = 0.752507
$ws.Cells.Item(17, 9).Value = 0.0642475070420447
$ws.Cells.Item(17, 10).Value = 0.0437690215344491
$ws.Cells.Item(17, 11).Value = 3
$ws.Cells.Item(17, 12).Value = 1
$ws.Cells.Item(17, 13).Value = 72.07364666666666
$ws.Cells.Item(17, 14).Value = 216.22094
$ws.Cells.Item(17, 15).Value = 0.4441129584835175
$ws.Cells.Item(17, 16).Value = 0.4455535891391496
$ws.Cells.Item(17, 17).Value = 27.11796181609667
$ws.Cells.Item(17, 18).Value = 162.70777089658
$ws.Cells.Item(17, 19).Value = 0.02853315042763309
$ws.Cells.Item(17, 20).Value = 0.01950144463778252

$ws.Cells.Item(18, 1).Value = "Resolving-Mac"
$ws.Cells.Item(18, 2).Value = "Spon2"
$ws.Cells.Item(18, 3).Value = "Itga4"
$ws.Cells.Item(18, 4).Value = "ECs"
$ws.Cells.Item(18, 5).Value = 1
$ws.Cells.Item(18, 6).Value = 0.3333333333333333
$ws.Cells.Item(18, 7).Value = 0.2467186666666667
$ws.Cells.Item(18, 8).Value = 0.740156
$ws.Cells.Item(18, 9).Value = 0.04212866929894484
$ws.Cells.Item(18, 10).Value = 0.04305063461582644
$ws.Cells.Item(18, 11).Value = 3
$ws.Cells.Item(18, 12).Value = 1
$ws.Cells.Item(18, 13).Value = 0.506715
$ws.Cells.Item(18, 14).Value = 1.520145
$ws.Cells.Item(18, 15).Value = 0.003122343715987576
$ws.Cells.Item(18, 16).Value = 0.003132472094339857
$ws.Cells.Item(18, 17).Value = 0.12501604918
$ws.Cells.Item(18, 18).Value = 1.12514444262
$ws.Cells.Item(18, 19).Value = 0.0001315401858484792
$ws.Cells.Item(18, 20).Value = 0.0001348549115776978

$ws.Cells.Item(19, 1).Value = "Resolving-Mac"
$ws.Cells.Item(19, 2).Value = "Spon2"
$ws.Cells.Item(19, 3).Value = "Itga4"
$ws.Cells.Item(19, 4).Value = "Inflammatory-Mac"
$ws.Cells.Item(19, 5).Value = 1
$ws.Cells.Item(19, 6).Value = 0.3333333333333333
$ws.Cells.Item(19, 7).Value = 0.2467186666666667
$ws.Cells.Item(19, 8).Value = 0.740156
$ws.Cells.Item(19, 9).Value = 0.04212866929894484
$ws.Cells.Item(19, 10).Value = 0.04305063461582644
$ws.Cells.Item(19, 11).Value = 3
$ws.Cells.Item(19, 12).Value = 1
$ws.Cells.Item(19, 13).Value = 88.13219433333332
$ws.Cells.Item(19, 14).Value = 264.396583
$ws.Cells.Item(19, 15).Value = 0.5430646480820168
$ws.Cells.Item(19, 16).Value = 0.5448262620252092
$ws.Cells.Item(19, 17).Value = 21.74385747632755
$ws.Cells.Item(19, 18).Value = 195.694717286948
$ws.Cells.Item(19, 19).Value = 0.02287859096699515
$ws.Cells.Item(19, 20).Value = 0.0234551163355538

$ws.Cells.Item(20, 1).Value = "Resolving-Mac"
$ws.Cells.Item(20, 2).Value = "Spon2"
$ws.Cells.Item(20, 3).Value = "Itga4"
$ws.Cells.Item(20, 4).Value = "MuSCs"
$ws.Cells.Item(20, 5).Value = 1
$ws.Cells.Item(20, 6).Value = 0.3333333333333333
$ws.Cells.Item(20, 7).Value = 0.2467186666666667
$ws.Cells.Item(20, 8).Value = 0.740156
$ws.Cells.Item(20, 9).Value = 0.04212866929894484
$ws.Cells.Item(20, 10).Value = 0.04305063461582644
$ws.Cells.Item(20, 11).Value = 2
$ws.Cells.Item(20, 12).Value = 1
$ws.Cells.Item(20, 13).Value = 1.5741895
$ws.Cells.Item(20, 14).Value = 3.148379
$ws.Cells.Item(20, 15).Value = 0.009700049718478087
$ws.Cells.Item(20, 16).Value = 0.006487676741301404
$ws.Cells.Item(20, 17).Value = 0.3883819345206667
$ws.Cells.Item(20, 18).Value = 2.330291607124
$ws.Cells.Item(20, 19).Value = 0.0004086501867730863
$ws.Cells.Item(20, 20).Value = 0.0002792986008953623

$ws.Cells.Item(21, 1).Value = "Resolving-Mac"
$ws.Cells.Item(21, 2).Value = "Spon2"
$ws.Cells.Item(21, 3).Value = "Itga4"
$ws.Cells.Item(21, 4).Value = "Resolving-Mac"
$ws.Cells.Item(21, 5).Value = 1
$ws.Cells.Item(21, 6).Value = 0.3333333333333333
$ws.Cells.Item(21, 7).Value = 0.2467186666666667
$ws.Cells.Item(21, 8).Value = 0.740156
$ws.Cells.Item(21, 9).Value = 0.04212866929894484
$ws.Cells.Item(21, 10).Value = 0.04305063461582644
$ws.Cells.Item(21, 11).Value = 3
$ws.Cells.Item(21, 12).Value = 1
$ws.Cells.Item(21, 13).Value = 72.07364666666666
$ws.Cells.Item(21, 14).Value = 216.22094
$ws.Cells.Item(21, 15).Value = 0.4441129584835175
$ws.Cells.Item(21, 16).Value = 0.4455535891391496
$ws.Cells.Item(21, 17).Value = 17.78191400740444
$ws.Cells.Item(21, 18).Value = 160.03722606664
$ws.Cells.Item(21, 19).Value = 0.01870988795932813
$ws.Cells.Item(21, 20).Value = 0.01918136476779958

